# Cryptos list update — regenerated from the scraper run on
# Thu Mar  9 08:53:45 UTC 2023 (GitHub Actions).
#
# Column D ("Price") values are plain text in this sheet (mixed
# thousands-dot / decimal-dot formatting from the source site), so any
# value that Excel would otherwise auto-parse as a number is written
# with a leading apostrophe to force text entry, exactly like typing
# it directly into the grid.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "21.704.39"
$ws.Range("E2").Value = "  -1.41%  "
# Row 3: Ethereum
$ws.Range("D3").Value = "1.534.95"
# Row 4: TetherUSD
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  +0.02%  "
# Row 5: USDC
$ws.Range("E5").Value = "  +0.06%  "
# Row 6: BNB
$ws.Range("D6").Value = "'289.05"
$ws.Range("E6").Value = "  +0.77%  "
# Row 7: XRP
$ws.Range("D7").Value = "'0.3944"
$ws.Range("E7").Value = "  +4.08%  "
# Row 8: Cardano
$ws.Range("D8").Value = "'0.3159"
$ws.Range("E8").Value = "  -2.56%  "
# Row 9: OKB
$ws.Range("D9").Value = "'42.19"
$ws.Range("E9").Value = "  +2.35%  "
# Row 10: Dogecoin
$ws.Range("D10").Value = "'0.07174"
$ws.Range("E10").Value = "  -1.86%  "
# Row 11: Polygon
$ws.Range("E11").Value = "  -6.83%  "
# Row 12: BinanceUSD
$ws.Range("D12").Value = "'1.000"
$ws.Range("E12").Value = "  +0.02%  "
# Row 13: Polkadot
$ws.Range("D13").Value = "'5.619"
$ws.Range("E13").Value = "  -1.87%  "
# Row 14: Solana
$ws.Range("D14").Value = "'18.45"
# Row 15: Chainlink
$ws.Range("D15").Value = "'6.598"
$ws.Range("E15").Value = "  -2.82%  "
# Row 16: WrappedEther
$ws.Range("D16").Value = "1.534.22"
$ws.Range("E16").Value = "  -2.14%  "
# Row 17: ShibaInu
$ws.Range("E17").Value = "  +0.83%  "
# Row 18: TRON
$ws.Range("D18").Value = "'0.06594"
$ws.Range("E18").Value = "  -0.34%  "
# Row 19: Litecoin
$ws.Range("D19").Value = "'83.05"
$ws.Range("E19").Value = "  -2.29%  "
# Row 20: Dai
$ws.Range("E20").Value = "  +0.22%  "
# Row 21: Uniswap
$ws.Range("D21").Value = "'6.115"
$ws.Range("E21").Value = "  -4.71%  "
# Row 22: Avalanche
$ws.Range("D22").Value = "'15.36"
$ws.Range("E22").Value = "  -3.82%  "
# Row 23: Cosmos
$ws.Range("E23").Value = "  -5.48%  "
# Row 24: Toncoin
$ws.Range("D24").Value = "'2.385"
$ws.Range("E24").Value = "  +4.32%  "
# Row 25: WrappedBTC
$ws.Range("D25").Value = "21.695.83"
$ws.Range("E25").Value = "  -1.52%  "
# Row 26: LidoDAOToken
$ws.Range("D26").Value = "'2.349"
$ws.Range("E26").Value = "  -7.07%  "
# Row 27: Monero
$ws.Range("D27").Value = "'146.42"
$ws.Range("E27").Value = "  -1.04%  "
# Row 28: EthereumClassic
$ws.Range("D28").Value = "'18.32"
$ws.Range("E28").Value = "  -2.56%  "
# Row 29: HuobiToken
$ws.Range("D29").Value = "'4.832"
$ws.Range("E29").Value = "  -0.54%  "
# Row 30: WrappedliquidstakedEther2.0
$ws.Range("D30").Value = "1.706.93"
$ws.Range("E30").Value = "  -2.04%  "
# Row 31: BitcoinCash
$ws.Range("D31").Value = "'116.85"
$ws.Range("E31").Value = "  -2.86%  "
# Row 32: ImmutableX
$ws.Range("D32").Value = "'0.9596"
$ws.Range("E32").Value = "  -13.43%  "
# Row 33: Filecoin
$ws.Range("D33").Value = "'5.848"
$ws.Range("E33").Value = "  -1.42%  "
# Row 34: Stellar
$ws.Range("D34").Value = "'0.08167"
$ws.Range("E34").Value = "  +0.24%  "
# Row 35: FraxShare
$ws.Range("D35").Value = "'8.634"
$ws.Range("E35").Value = "  -6.31%  "
# Row 36: Hedera
$ws.Range("D36").Value = "'0.06064"
$ws.Range("E36").Value = "  -2.01%  "
# Row 37: InternetComputer(DFINITY)
$ws.Range("E37").Value = "  -2.90%  "
# Row 38: VeChain
$ws.Range("D38").Value = "'0.02195"
# Row 39: WEMIXTOKEN
$ws.Range("B39").Value = "WEMIXTOKEN"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").Value = "'1.437"
$ws.Range("E39").Value = "  -12.80%  "
# Row 40: Algorand
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").Value = "'0.2018"
$ws.Range("E40").Value = "  -4.51%  "
# Row 41: TrustWalletToken
$ws.Range("D41").Value = "'1.178"
$ws.Range("E41").Value = "  -3.44%  "
# Row 42: Frax
$ws.Range("E42").Value = "  +0.10%  "
# Row 43: Aptos
$ws.Range("D43").Value = "'10.69"
$ws.Range("E43").Value = "  -1.74%  "
# Row 44: TheSandbox
$ws.Range("D44").Value = "'0.5717"
$ws.Range("E44").Value = "  -3.64%  "
# Row 45: PancakeSwap
$ws.Range("B45").Value = "PancakeSwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D45").Value = "'3.732"
$ws.Range("E45").Value = "  +0.32%  "
# Row 46: EnergySwap
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").Value = "'13.02"
$ws.Range("E46").Value = "  -3.40%  "
# Row 47: Decentraland
$ws.Range("D47").Value = "'0.5465"
$ws.Range("E47").Value = "  -4.67%  "
# Row 48: EOS
$ws.Range("D48").Value = "'1.156"
$ws.Range("E48").Value = "  -0.11%  "
# Row 49: Quant
$ws.Range("D49").Value = "'115.86"
$ws.Range("E49").Value = "  -3.29%  "
# Row 50: NEARProtocol
$ws.Range("D50").Value = "'1.861"
$ws.Range("E50").Value = "  -3.74%  "
# Row 51: Cronos
$ws.Range("D51").Value = "'0.06692"
$ws.Range("E51").Value = "  -2.74%  "
